# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E27) listed the account-statement periods
# in descending order (1805 -> 1706). This update refreshes the database so
# the periods are listed in ascending order (1706 -> 1805), i.e. each row's
# period value is replaced with the next one in the new, ascending sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periods = @("1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805")

$row = 16
foreach ($periodo in $periods) {
    $ws.Range("E$row").Value = $periodo
    $row++
}
